$d = $word.ActiveDocument

$d.Content.Find.Execute("452÷8=56, 4", $true, $false, $false, $false, $false, $true, 1, $false, "844÷6=140, 4", 2) | Out-Null
$d.Content.Find.Execute("105÷6=17, 3", $true, $false, $false, $false, $false, $true, 1, $false, "865÷4=216, 1", 2) | Out-Null
$d.Content.Find.Execute("148÷8=18, 4", $true, $false, $false, $false, $false, $true, 1, $false, "368÷3=122, 2", 2) | Out-Null
$d.Content.Find.Execute("885÷8=110, 5", $true, $false, $false, $false, $false, $true, 1, $false, "764÷5=152, 4", 2) | Out-Null
$d.Content.Find.Execute("643÷2=321, 1", $true, $false, $false, $false, $false, $true, 1, $false, "799÷7=114, 1", 2) | Out-Null
$d.Content.Find.Execute("112÷7=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "156÷6=26, 0", 2) | Out-Null
$d.Content.Find.Execute("584÷6=97, 2", $true, $false, $false, $false, $false, $true, 1, $false, "875÷4=218, 3", 2) | Out-Null
$d.Content.Find.Execute("309÷6=51, 3", $true, $false, $false, $false, $false, $true, 1, $false, "520÷2=260, 0", 2) | Out-Null
$d.Content.Find.Execute("262÷6=43, 4", $true, $false, $false, $false, $false, $true, 1, $false, "198÷8=24, 6", 2) | Out-Null
$d.Content.Find.Execute("583÷5=116, 3", $true, $false, $false, $false, $false, $true, 1, $false, "876÷5=175, 1", 2) | Out-Null
$d.Content.Find.Execute("857÷6=142, 5", $true, $false, $false, $false, $false, $true, 1, $false, "641÷4=160, 1", 2) | Out-Null
$d.Content.Find.Execute("138÷4=34, 2", $true, $false, $false, $false, $false, $true, 1, $false, "552÷2=276, 0", 2) | Out-Null
$d.Content.Find.Execute("138÷3=46, 0", $true, $false, $false, $false, $false, $true, 1, $false, "191÷6=31, 5", 2) | Out-Null
$d.Content.Find.Execute("374÷4=93, 2", $true, $false, $false, $false, $false, $true, 1, $false, "637÷4=159, 1", 2) | Out-Null
$d.Content.Find.Execute("215÷5=43, 0", $true, $false, $false, $false, $false, $true, 1, $false, "432÷2=216, 0", 2) | Out-Null
$d.Content.Find.Execute("321÷4=80, 1", $true, $false, $false, $false, $false, $true, 1, $false, "477÷5=95, 2", 2) | Out-Null
$d.Content.Find.Execute("136÷3=45, 1", $true, $false, $false, $false, $false, $true, 1, $false, "402÷4=100, 2", 2) | Out-Null
$d.Content.Find.Execute("157÷3=52, 1", $true, $false, $false, $false, $false, $true, 1, $false, "138÷4=34, 2", 2) | Out-Null
$d.Content.Find.Execute("738÷7=105, 3", $true, $false, $false, $false, $false, $true, 1, $false, "483÷8=60, 3", 2) | Out-Null
$d.Content.Find.Execute("387÷6=64, 3", $true, $false, $false, $false, $false, $true, 1, $false, "315÷3=105, 0", 2) | Out-Null
$d.Content.Find.Execute("361÷7=51, 4", $true, $false, $false, $false, $false, $true, 1, $false, "796÷5=159, 1", 2) | Out-Null
$d.Content.Find.Execute("264÷5=52, 4", $true, $false, $false, $false, $false, $true, 1, $false, "413÷8=51, 5", 2) | Out-Null
$d.Content.Find.Execute("422÷4=105, 2", $true, $false, $false, $false, $false, $true, 1, $false, "189÷2=94, 1", 2) | Out-Null
$d.Content.Find.Execute("855÷6=142, 3", $true, $false, $false, $false, $false, $true, 1, $false, "888÷5=177, 3", 2) | Out-Null
$d.Content.Find.Execute("311÷3=103, 2", $true, $false, $false, $false, $false, $true, 1, $false, "637÷5=127, 2", 2) | Out-Null
